$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Furniture sheet updates
# ---------------------------------------------------------------------------
$wsFurniture = $wb.Worksheets.Item("Furniture")
[void]$wsFurniture.Activate()

# Raw-material 1st lot (Glue) cost increased
$wsFurniture.Range("F10").Value = 18200

# Insert two new blank rows right after row 33 (pushes old rows 34-35 down to 36-37)
$wsFurniture.Rows("34:35").Insert()

# Fill in the newly used row (old row 33, still row 33 after insert) with the
# extra glue purchase info
$wsFurniture.Range("F33").Value = 1750
$wsFurniture.Range("G33").Value = "Self"

# Restore the on-screen selection to match where the user ended up
[void]$wsFurniture.Range("K14").Select()

# ---------------------------------------------------------------------------
# 2. Comment update on Furniture!F10
# ---------------------------------------------------------------------------
$comment = $wsFurniture.Range("F10").Comment
[void]$comment.Text("sounak nandi:" + "`n" + "3 tubs 20kg" + "`n" + "5kg special glue" + "`n")

# ---------------------------------------------------------------------------
# 3. Index sheet updates
# ---------------------------------------------------------------------------
$wsIndex = $wb.Worksheets.Item("Index")
[void]$wsIndex.Activate()

# Grill advance amount updated
$wsIndex.Range("H13").Value = 38050

# Restore the on-screen selection to match where the user ended up
[void]$wsIndex.Range("H17:I17").Select()
